$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Suite")

# Update C4 value from "N" to "Y"
$ws.Range("C4").Value = "Y"

# Move active selection from C2 to C3
$ws.Range("C3").Select()
